$d = $word.ActiveDocument

# The document has four list-item paragraphs whose whole text is
# "Feature N [pasca bayar]" / "Feature N [prabayar]". Each one must be
# split into two runs: "Menu" and " N [...]" (note the leading space on
# the second run, preserved via xml:space="preserve").
#
# We rebuild each paragraph's run content (but not its paragraph mark,
# so w:pPr / numbering stays untouched) via InsertXML, which lets us
# drop in exactly the two <w:r> runs we want without leaving stray
# empty <w:rPr/> residue behind.

foreach ($p in $d.Paragraphs) {
    $full = $p.Range.Text
    if ($full -match "^Feature( [^\r]*)\r?$") {
        $rest = $matches[1]
        $r = $p.Range
        # exclude the trailing paragraph mark from the replaced range
        $contentRange = $d.Range($r.Start, $r.End - 1)

        $escRest = $rest -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"

        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Menu</w:t></w:r><w:r><w:t xml:space="preserve">' + $escRest + '</w:t></w:r></w:p>'
        $contentRange.InsertXML($xml)
    }
}
